# Alteração nos rótulos da tabela para já transformar a primeira linha em
# cabeçalho automaticamente no Power BI: prefixa os cabeçalhos de ano com
# "Ano" (ou "Intervalo" na planilha de potência incremental).

$wb = $excel.ActiveWorkbook

# Planilhas cujo cabeçalho (B1:E1) recebe o prefixo "Ano "
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range("$col" + "1")
        $cell.Value = "Ano " + $cell.Value2
    }
}

# Planilha de potência incremental recebe o prefixo "Intervalo "
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws.Range("$col" + "1")
    $cell.Value = "Intervalo " + $cell.Value2
}

# Planilha de custo total só tem a coluna B com o ano
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$cell = $ws.Range("B1")
$cell.Value = "Ano " + $cell.Value2
